# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" worksheet right after "2021-Q4" (before "总计"),
#   populated with the fund-holdings breakdown table.
# - Insert a new summary row for "2022-Q1" at the top of the "总计" sheet's
#   data (pushing the existing rows down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: apply the bold / centered / thin-border "header-ish" look used
# throughout this workbook (column-A index cells + row-1 headers).
# ---------------------------------------------------------------------
function Set-HeaderStyle($cell) {
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous
}

# ---------------------------------------------------------------------
# Helper: write a value that must be stored as literal TEXT even though
# it looks numeric (fund codes with leading zeros, numeric-looking
# ratios/values that the source data keeps as strings). Forces the
# string type, then clears the leftover "@" number-format styling so
# the cell ends up with no special style - matching plain data cells.
# ---------------------------------------------------------------------
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# =======================================================================
# 1. Insert the new "2022-Q1" sheet after "2021-Q4" and before "总计"
# =======================================================================
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $afterSheet)
$q1.Name = "2022-Q1"

# ---- header row ----
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q1.Cells.Item(1, 2 + $i)
    $cell.Value = $headers[$i]
    Set-HeaderStyle $cell
}

# ---- data rows ----
# columns: idx, code, name, scale, stockPos, posRatio, marketValue, posRank
$rows = @(
    @(0, "004856", "广发中证全指建筑材料指数A", "13.72", "94.61", "2.71", "0.3718", 8),
    @(1, "004857", "广发中证全指建筑材料指数C", "6.05",  "94.61", "2.71", "0.1640", 8),
    @(2, "159745", "国泰中证全指建筑材料交易型开放式指数证券投资基金", "3.76", "98.37", "2.79", "0.1049", 9),
    @(3, "970020", "信达价值精选一年持有期灵活配置混合A", "0.64", "56.02", "5.18", "0.0332", 4),
    @(4, "970021", "信达价值精选一年持有期灵活配置混合B", "0.53", "56.02", "5.18", "0.0275", 4),
    @(5, "516750", "富国中证全指建筑材料ETF", "0.47", "98.22", "2.82", "0.0133", 8)
)

$r = 2
foreach ($row in $rows) {
    $idxCell = $q1.Cells.Item($r, 1)
    $idxCell.Value = $row[0]
    Set-HeaderStyle $idxCell

    Set-TextValue $q1.Cells.Item($r, 2) $row[1]
    Set-TextValue $q1.Cells.Item($r, 3) $row[2]
    Set-TextValue $q1.Cells.Item($r, 4) $row[3]
    Set-TextValue $q1.Cells.Item($r, 5) $row[4]
    Set-TextValue $q1.Cells.Item($r, 6) $row[5]
    Set-TextValue $q1.Cells.Item($r, 7) $row[6]
    $q1.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

# =======================================================================
# 2. Update the "总计" sheet: insert a new "2022-Q1" row at the top of the
#    data (row 2), shifting the existing history rows down by one.
# =======================================================================
$total = $wb.Worksheets.Item("总计")

# Read the existing data rows (old rows 2..5) before they get overwritten.
$oldRowCount = 4
$existing = @()
for ($i = 0; $i -lt $oldRowCount; $i++) {
    $srcRow = 2 + $i
    $existing += , @($total.Cells.Item($srcRow, 2).Value2, $total.Cells.Item($srcRow, 3).Value2, $total.Cells.Item($srcRow, 4).Value2)
}

# New full data set for the "总计" sheet (2022-Q1 on top, then history).
$totalRows = @(
    @("2022-Q1", 6, 0.71)
)
$totalRows += $existing

$r = 2
foreach ($row in $totalRows) {
    $total.Cells.Item($r, 1).Value = $r - 2
    Set-HeaderStyle $total.Cells.Item($r, 1)
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r = $r + 1
}

Write-Output "2022-Q1 sheet added; 总计 updated"
